$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that needs to move
# from 45182 (2023-09-13) to 45184 (2023-09-15) for every data row
# (rows 2 through 246).
$lastRow = 246
$range = $ws.Range("C2:C$lastRow")
$range.Value2 = 45184
